$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IXTK102N65X2")
$ws.Range("K5").Value = 0.16948353896707699
$ws.Range("K6").Value = 0.090038130076260103
$ws.Range("K7").Select()
